$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 20241024
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 6

# Row 6
$ws.Range("D6").Value = 29

# Row 7
$ws.Range("A7").Value = 20241120
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 5

# Row 9
$ws.Range("A9").Value = 20241121
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 5

# Row 11
$ws.Range("A11").Value = 20250123
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 6

# Row 12
$ws.Range("B12").Value = 15

# Update selection to match final cursor position
$ws.Range("B12").Select()
